$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rule table "SimpleRules String TestSimple(String param1)" added below
# the existing "Hello" table (rows 13-15), mirroring the look of a generated
# business-rule spreadsheet table: a merged title row, a bold header row
# (param name / RETURN) and one data row (step1 / run).
# ---------------------------------------------------------------------------

# --- Row 13: merged title row -------------------------------------------------
$title = $ws.Range("B13:C13")
$title.Merge()
$title.Value = "SimpleRules String TestSimple(String param1)"
$title.Style = "Normal"
$title.Font.Name = "Franklin Gothik Book"
$title.Font.Size = 10
$title.Font.Bold = $false
$title.Font.Underline = $false
$title.Font.Color = 65536
$title.Interior.Pattern = 1
$title.Interior.Color = 16711679
$title.Borders.Item(8).LineStyle = 1
$title.Borders.Item(8).Weight = 2
$title.Borders.Item(8).Color = 65536
$title.Borders.Item(9).LineStyle = 1
$title.Borders.Item(9).Weight = 2
$title.Borders.Item(9).Color = 65536
$title.HorizontalAlignment = -4108

# --- Row 14: header row (param1 | RETURN) -------------------------------------
$paramHeader = $ws.Range("B14")
$paramHeader.Value = "param1"
$paramHeader.Style = "Normal"
$paramHeader.Font.Name = "Franklin Gothik Book"
$paramHeader.Font.Size = 10
$paramHeader.Font.Bold = $true
$paramHeader.Font.Underline = $false
$paramHeader.Font.Color = 0
$paramHeader.Interior.Pattern = 1
$paramHeader.Interior.Color = 10921638
$paramHeader.Borders.Item(9).LineStyle = 1
$paramHeader.Borders.Item(9).Weight = 2
$paramHeader.Borders.Item(9).Color = 65536
$paramHeader.HorizontalAlignment = -4108

$returnHeader = $ws.Range("C14")
$returnHeader.Value = "RETURN"
$returnHeader.Style = "Normal"
$returnHeader.Font.Name = "Franklin Gothik Book"
$returnHeader.Font.Size = 10
$returnHeader.Font.Bold = $true
$returnHeader.Font.Underline = $false
$returnHeader.Font.Color = 0
$returnHeader.Interior.Pattern = 1
$returnHeader.Interior.Color = 4641530
$returnHeader.Borders.Item(9).LineStyle = 1
$returnHeader.Borders.Item(9).Weight = 4
$returnHeader.Borders.Item(9).Color = 1952255
$returnHeader.HorizontalAlignment = -4108

# --- Row 15: data row (step1 | run) ------------------------------------------
$step1 = $ws.Range("B15")
$step1.Value = "step1"
$step1.Style = "Normal"
$step1.Font.Name = "Franklin Gothik Book"
$step1.Font.Size = 10
$step1.Font.Bold = $false
$step1.Font.Underline = $false
$step1.Font.Color = 0
$step1.Interior.Pattern = 1
$step1.Interior.Color = 16711679
$step1.Borders.Item(7).LineStyle = 1
$step1.Borders.Item(7).Weight = 2
$step1.Borders.Item(7).Color = 14540253
$step1.Borders.Item(10).LineStyle = 1
$step1.Borders.Item(10).Weight = 2
$step1.Borders.Item(10).Color = 14540253
$step1.Borders.Item(8).LineStyle = 1
$step1.Borders.Item(8).Weight = 2
$step1.Borders.Item(8).Color = 14540253
$step1.Borders.Item(9).LineStyle = 1
$step1.Borders.Item(9).Weight = 2
$step1.Borders.Item(9).Color = 65536
$step1.HorizontalAlignment = -4108

$run = $ws.Range("C15")
$run.Value = "run"
$run.Style = "Normal"
$run.Font.Name = "Franklin Gothik Book"
$run.Font.Size = 10
$run.Font.Bold = $false
$run.Font.Underline = $false
$run.Font.Color = 0
$run.Interior.Pattern = 1
$run.Interior.Color = 14277081
$run.Borders.Item(7).LineStyle = 1
$run.Borders.Item(7).Weight = 2
$run.Borders.Item(7).Color = 14540253
$run.Borders.Item(10).LineStyle = 1
$run.Borders.Item(10).Weight = 2
$run.Borders.Item(10).Color = 14540253
$run.Borders.Item(8).LineStyle = 1
$run.Borders.Item(8).Weight = 2
$run.Borders.Item(8).Color = 14540253
$run.Borders.Item(9).LineStyle = 1
$run.Borders.Item(9).Weight = 4
$run.Borders.Item(9).Color = 1952255
$run.HorizontalAlignment = -4108

Write-Output "TestSimple table added"
